$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "64 x 18" + [char]11 + "  1    8" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "4|    |"
$t.Cell(1,2).Range.Text = "29 x 44" + [char]11 + "  4    4" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "9|    |"
$t.Cell(1,3).Range.Text = "89 x 60" + [char]11 + "  6    0" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "9|    |"
$t.Cell(2,1).Range.Text = "55 x 34" + [char]11 + "  3    4" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "5|    |"
$t.Cell(2,2).Range.Text = "25 x 92" + [char]11 + "  9    2" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "5|    |"
$t.Cell(2,3).Range.Text = "15 x 22" + [char]11 + "  2    2" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "5|    |"
$t.Cell(3,1).Range.Text = "57 x 86" + [char]11 + "  8    6" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "7|    |"
$t.Cell(3,2).Range.Text = "48 x 65" + [char]11 + "  6    5" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "8|    |"
$t.Cell(3,3).Range.Text = "90 x 61" + [char]11 + "  6    1" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "0|    |"
$t.Cell(4,1).Range.Text = "89 x 83" + [char]11 + "  8    3" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "9|    |"
$t.Cell(4,2).Range.Text = "68 x 97" + [char]11 + "  9    7" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "8|    |"
$t.Cell(4,3).Range.Text = "83 x 44" + [char]11 + "  4    4" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "3|    |"
$t.Cell(5,1).Range.Text = "41 x 44" + [char]11 + "  4    4" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "1|    |"
$t.Cell(5,2).Range.Text = "68 x 82" + [char]11 + "  8    2" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "8|    |"
$t.Cell(5,3).Range.Text = "94 x 42" + [char]11 + "  4    2" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "4|    |"
